# Convert the Word FIELD "{ m:self.name }" (fldChar begin/instrText.../fldChar end)
# in the second paragraph into plain literal text runs:
#   "{"  "m"  ":"  "self"(colored)  ".name}"
#
# Word's field-code runs are atomic from a COM editing point of view (any
# deletion that touches the field removes the whole field), so the field is
# deleted in one shot and the equivalent plain-text runs are retyped in its
# place. A temporary bookmark is dropped at each run boundary (and removed
# right after) purely to stop same-formatting neighboring runs from being
# coalesced into a single run when the document is saved, so the run
# boundaries shown in the target markup are preserved.

$d = $word.ActiveDocument

# Locate the field that holds " m:self.name ".
$targetField = $null
foreach ($fld in $d.Fields) {
    if ($fld.Code.Text -match "self") {
        $targetField = $fld
    }
}

# The field's begin fldChar sits exactly one position before its Code range
# starts; once the field is removed that position becomes the insertion
# point for the replacement text.
$pos = $targetField.Code.Start - 1

# Delete the whole field (begin/instrText.../end) - this collapses the
# paragraph down to just its paragraph mark.
$targetField.Delete()

# --- "{" -----------------------------------------------------------------
$run = $d.Range($pos, $pos)
$run.InsertAfter("{")
$pos = $pos + 1

$boundary = $d.Range($pos, $pos)
$d.Bookmarks.Add("zzzM2DocSplit1", $boundary)

# --- "m" -------------------------------------------------------------------
$run = $d.Range($pos, $pos)
$run.InsertAfter("m")
$pos = $pos + 1

$boundary = $d.Range($pos, $pos)
$d.Bookmarks.Add("zzzM2DocSplit2", $boundary)

# --- ":" -------------------------------------------------------------------
$run = $d.Range($pos, $pos)
$run.InsertAfter(":")
$pos = $pos + 1

# --- "self" (keeps the accent6 theme color the field code had) -----------
$run = $d.Range($pos, $pos)
$run.InsertAfter("self")
$selfRange = $d.Range($pos, $pos + 4)
$selfRange.Font.TextColor.ObjectThemeColor = 9
$pos = $pos + 4

# --- ".name}" --------------------------------------------------------------
$run = $d.Range($pos, $pos)
$run.InsertAfter(".name}")
$pos = $pos + 6

# Remove the temporary boundary bookmarks - they only existed to keep the
# runs above from being merged back together on save.
$d.Bookmarks("zzzM2DocSplit1").Delete()
$d.Bookmarks("zzzM2DocSplit2").Delete()
